$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column AA with 2023 data, mirroring the format of column Z.
$ws.Range("AA4").Value = 2023
$ws.Range("AA5").Value = 44.2
$ws.Range("AA6").Value = 50.4
$ws.Range("AA7").Value = 40.6
$ws.Range("AA8").Value = 57.2
$ws.Range("AA9").Value = 31
$ws.Range("AA10").Value = 49.7
$ws.Range("AA11").Value = 51
$ws.Range("AA12").Value = 29.4
$ws.Range("AA13").Value = 29.9
$ws.Range("AA14").Value = 56.3
$ws.Range("AA15").Value = 62.5
$ws.Range("AA16").Value = 34.9

# Copy the style/formatting from column Z (the prior last column) to column AA.
$ws.Range("Z4:Z16").Copy()
$ws.Range("AA4:AA16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Reset the selection/active cell and scroll position so the view matches.
$ws.Range("A1").Select()
